$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Worksheet 1
$ws1.Cells.Item(2, 6).Value = 1605
$ws1.Cells.Item(3, 6).Value = 3343
$ws1.Cells.Item(4, 6).Value = 33
$ws1.Cells.Item(5, 6).Value = 766
$ws1.Cells.Item(6, 6).Value = 2389
$ws1.Cells.Item(8, 6).Value = 431
$ws1.Cells.Item(10, 6).Value = 156
$ws1.Cells.Item(11, 6).Value = 379
$ws1.Cells.Item(12, 6).Value = 1122
$ws1.Cells.Item(13, 6).Value = 476
$ws1.Cells.Item(15, 6).Value = 98
$ws1.Cells.Item(16, 6).Value = 288
$ws1.Cells.Item(17, 6).Value = 4979
$ws1.Cells.Item(18, 6).Value = 31
$ws1.Cells.Item(19, 6).Value = 1401
$ws1.Cells.Item(20, 6).Value = 3662
$ws1.Cells.Item(21, 6).Value = 196
$ws1.Cells.Item(22, 6).Value = 227
$ws1.Cells.Item(23, 6).Value = 4000
$ws1.Cells.Item(24, 6).Value = 5348
$ws1.Cells.Item(27, 6).Value = 586
$ws1.Cells.Item(28, 6).Value = 3428
$ws1.Cells.Item(29, 6).Value = 404
$ws1.Cells.Item(31, 6).Value = 154
$ws1.Cells.Item(33, 6).Value = 910
$ws1.Cells.Item(34, 6).Value = 1240
$ws1.Cells.Item(36, 6).Value = 70
$ws1.Cells.Item(37, 6).Value = 1467
$ws1.Cells.Item(39, 6).Value = 1456
$ws1.Cells.Item(40, 6).Value = 50
$ws1.Cells.Item(41, 6).Value = 942
$ws1.Cells.Item(42, 6).Value = 937
$ws1.Cells.Item(43, 6).Value = 537
$ws1.Cells.Item(44, 6).Value = 67
$ws1.Cells.Item(45, 6).Value = 2504
$ws1.Cells.Item(46, 6).Value = 95
$ws1.Cells.Item(47, 6).Value = 196
$ws1.Cells.Item(49, 6).Value = 3777

# Worksheet 2
$ws2.Cells.Item(6, 6).Value = 1046
$ws2.Cells.Item(10, 6).Value = 10

# Worksheet 3
$ws3.Cells.Item(2, 6).Value = 2847

# Worksheet 4
$ws4.Cells.Item(2, 6).Value = 2847
$ws4.Cells.Item(3, 6).Value = 1605
$ws4.Cells.Item(4, 6).Value = 3343
$ws4.Cells.Item(5, 6).Value = 33
$ws4.Cells.Item(6, 6).Value = 766
$ws4.Cells.Item(7, 6).Value = 2389
$ws4.Cells.Item(9, 6).Value = 431
$ws4.Cells.Item(11, 6).Value = 1046
$ws4.Cells.Item(12, 6).Value = 156
$ws4.Cells.Item(13, 6).Value = 379
$ws4.Cells.Item(14, 6).Value = 1122
$ws4.Cells.Item(15, 6).Value = 476
$ws4.Cells.Item(17, 6).Value = 98
$ws4.Cells.Item(18, 6).Value = 288
$ws4.Cells.Item(19, 6).Value = 4979
$ws4.Cells.Item(20, 6).Value = 1401
$ws4.Cells.Item(21, 6).Value = 4000
$ws4.Cells.Item(22, 6).Value = 5348
$ws4.Cells.Item(25, 6).Value = 586
$ws4.Cells.Item(26, 6).Value = 3428
$ws4.Cells.Item(27, 6).Value = 404
$ws4.Cells.Item(29, 6).Value = 154
$ws4.Cells.Item(31, 6).Value = 910
$ws4.Cells.Item(32, 6).Value = 1240
$ws4.Cells.Item(34, 6).Value = 70
$ws4.Cells.Item(35, 6).Value = 1467
$ws4.Cells.Item(37, 6).Value = 1456
$ws4.Cells.Item(39, 6).Value = 942
$ws4.Cells.Item(41, 6).Value = 537
$ws4.Cells.Item(43, 6).Value = 67
$ws4.Cells.Item(45, 6).Value = 2504
$ws4.Cells.Item(46, 6).Value = 95
$ws4.Cells.Item(47, 6).Value = 196
$ws4.Cells.Item(49, 6).Value = 3777
